# Investimentos.xlsx - "Mudanças para melhor experiência do usuário"
#
# - Remove the two bottom rows (vale3 / beef3), shrinking the table to a
#   single data row (petr4).
# - petr4's quoted price (Valor) moved from 36.51 to 38.5, so the
#   price-derived ratios in that row (EY, EY2, P/L, P/VP, EV/EBITDA) are
#   refreshed to match.
# - Conditional-formatting ranges that used to span the 3 data rows now
#   only cover the single remaining row, and the per-row highlight rules
#   that belonged to the deleted rows are removed. The thresholds for the
#   rules tied to petr4's own price (VPA / Teto 9% / Valor Justo columns)
#   are updated to the new 38.5 price point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Drop the vale3 (row 3) and beef3 (row 4) rows entirely.
# ---------------------------------------------------------------------
$ws.Rows("3:4").Delete()

# ---------------------------------------------------------------------
# 2. Refresh petr4's (row 2) values for the new quoted price of 38.5.
# ---------------------------------------------------------------------
$ws.Range("B2").Value = 38.5
$ws.Range("I2").Value = 15.68831168831169
$ws.Range("J2").Value = 81.16883116883116
$ws.Range("K2").Value = 6.38
$ws.Range("L2").Value = 1.23
$ws.Range("P2").Value = 3.23

# ---------------------------------------------------------------------
# 3. Shrink the conditional-formatting ranges that used to span rows
#    2-4 down to just row 2 (the only remaining data row).
# ---------------------------------------------------------------------
$shrinkCols = @("F","Q","N","R","O","K","L","S","M","P","I","J")
foreach ($col in $shrinkCols) {
    $oldRange = $ws.Range($col + "2:" + $col + "4")
    $fcs = $oldRange.FormatConditions
    $newRange = $ws.Range($col + "2")
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 4. Remove the per-cell highlight rules that belonged to the now-gone
#    vale3 / beef3 rows (VPA / Teto 9% / Valor Justo columns).
# ---------------------------------------------------------------------
$removedRefs = @("C3", "D3", "E3", "C4", "D4", "E4")
foreach ($ref in $removedRefs) {
    $fcs = $ws.Range($ref).FormatConditions
    for ($i = $fcs.Count; $i -ge 1; $i--) {
        $fcs.Item($i).Delete()
    }
}

# ---------------------------------------------------------------------
# 5. petr4's own VPA / Teto 9% / Valor Justo rules compare against its
#    own price -- move the thresholds from 36.51/36.509 to 38.5/38.499.
# ---------------------------------------------------------------------
$ownPriceRefs = @("C2", "D2", "E2")
foreach ($ref in $ownPriceRefs) {
    $fcs = $ws.Range($ref).FormatConditions
    $fcs.Item(1).Formula1 = "=38.5"
    $fcs.Item(2).Formula1 = "=0.001"
    $fcs.Item(2).Formula2 = "=38.499"
}

Write-Output "edit applied"
